$d = $word.ActiveDocument
$tbl = $d.Tables(1)

# --- Simple cell text replacements (rows 1-6, by position so it is
#     unambiguous even though several cells share the same old text) ---
$tbl.Rows(1).Cells(1).Range.Text  = "0M"      # was 100
$tbl.Rows(2).Cells(1).Range.Text  = "0M"      # was 0.03
$tbl.Rows(3).Cells(1).Range.Text  = "0M"      # was 2447
$tbl.Rows(4).Cells(1).Range.Text  = "562"     # was 3
$tbl.Rows(5).Cells(1).Range.Text  = "0.00002" # was 0.00003
$tbl.Rows(6).Cells(1).Range.Text  = "0.00011" # was 0.00005

# --- Row 9 (text "0.00003") is removed entirely ---
$tbl.Rows(9).Delete()

# After the delete, the row that used to be row 12 ("0.00012") is now row 11.
$tbl.Rows(11).Cells(1).Range.Text = "0.00005" # was 0.00012

# Insert a brand-new row right after it (i.e. before the old row 13,
# now row 12) and populate it.
$newRow = $tbl.Rows.Add($tbl.Rows(12))
$tbl.Rows(12).Cells(1).Range.Text = "0.02766"

# --- Collapse the three multi-run summary rows (originally rows 44-46,
#     unaffected in index since one row was deleted and one inserted
#     earlier) down to a single value each ---
$tbl.Rows(44).Cells(1).Range.Text = "100"
$tbl.Rows(45).Cells(1).Range.Text = "0.03"
$tbl.Rows(46).Cells(1).Range.Text = "2447"
